$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'267.06"
$ws.Cells.Item(2, 7).Value = "'21"
$ws.Cells.Item(3, 4).Value = "'22.84"
$ws.Cells.Item(3, 7).Value = "'21"
$ws.Cells.Item(4, 4).Value = "'6.333"
$ws.Cells.Item(4, 7).Value = "'21"
$ws.Cells.Item(5, 4).Value = "'0.06199"
$ws.Cells.Item(5, 7).Value = "'21"
$ws.Cells.Item(6, 4).Value = "'3.590"
$ws.Cells.Item(6, 7).Value = "'21"
$ws.Cells.Item(7, 4).Value = "'6.665"
$ws.Cells.Item(7, 7).Value = "'21"
$ws.Cells.Item(8, 4).Value = "'1.389"
$ws.Cells.Item(8, 7).Value = "'21"
$ws.Cells.Item(9, 4).Value = "'0.8282"
$ws.Cells.Item(9, 7).Value = "'21"
$ws.Cells.Item(10, 7).Value = "'21"
$ws.Cells.Item(11, 4).Value = "'0.1608"
$ws.Cells.Item(11, 7).Value = "'21"
$ws.Cells.Item(12, 4).Value = "'0.08191"
$ws.Cells.Item(12, 7).Value = "'21"
$ws.Cells.Item(13, 4).Value = "'0.03410"
$ws.Cells.Item(13, 7).Value = "'21"
$ws.Cells.Item(14, 4).Value = "'0.03159"
$ws.Cells.Item(14, 7).Value = "'21"
$ws.Cells.Item(15, 4).Value = "'0.09276"
$ws.Cells.Item(15, 7).Value = "'21"
$ws.Cells.Item(16, 4).Value = "'3.900"
$ws.Cells.Item(16, 7).Value = "'21"
$ws.Cells.Item(17, 7).Value = "'21"
$ws.Cells.Item(18, 4).Value = "'0.04844"
$ws.Cells.Item(18, 7).Value = "'21"
$ws.Cells.Item(19, 4).Value = "'0.006218"
$ws.Cells.Item(19, 7).Value = "'21"
$ws.Cells.Item(20, 4).Value = "'0.005383"
$ws.Cells.Item(20, 7).Value = "'21"
$ws.Cells.Item(21, 2).Value = "UpBots"
$ws.Cells.Item(21, 3).Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
$ws.Cells.Item(21, 4).Value = "'0.007502"
$ws.Cells.Item(21, 5).Value = "20UpBotsUBXT"
$ws.Cells.Item(21, 7).Value = "'21"
$ws.Cells.Item(22, 2).Value = "BitKan"
$ws.Cells.Item(22, 3).Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Cells.Item(22, 4).Value = "'0.001091"
$ws.Cells.Item(22, 5).Value = "21BitKanKAN"
$ws.Cells.Item(22, 7).Value = "'21"
$ws.Cells.Item(23, 2).Value = "NitroEx"
$ws.Cells.Item(23, 3).Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Cells.Item(23, 4).Value = "'0.0001502"
$ws.Cells.Item(23, 5).Value = "22NitroExNTX"
$ws.Cells.Item(23, 7).Value = "'21"
$ws.Cells.Item(24, 2).Value = "LEO"
$ws.Cells.Item(24, 3).Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Cells.Item(24, 4).Value = "'3.758"
$ws.Cells.Item(24, 5).Value = "23LEOLEO"
$ws.Cells.Item(24, 7).Value = "'21"
$ws.Cells.Item(25, 2).Value = "BTSEToken"
$ws.Cells.Item(25, 3).Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Cells.Item(25, 4).Value = "'2.369"
$ws.Cells.Item(25, 5).Value = "24BTSETokenBTSE"
$ws.Cells.Item(25, 7).Value = "'21"
$ws.Cells.Item(26, 2).Value = "BitpandaEcosystemToken"
$ws.Cells.Item(26, 3).Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Cells.Item(26, 4).Value = "'0.3350"
$ws.Cells.Item(26, 5).Value = "25BitpandaEcosystemTokenBEST"
$ws.Cells.Item(26, 7).Value = "'21"
$ws.Cells.Item(27, 2).Value = "ProBitToken"
$ws.Cells.Item(27, 3).Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Cells.Item(27, 4).Value = "'0.1214"
$ws.Cells.Item(27, 5).Value = "26ProBitTokenPROB"
$ws.Cells.Item(27, 7).Value = "'21"
$ws.Cells.Item(28, 7).Value = "'21"
$ws.Cells.Item(29, 7).Value = "'21"
$ws.Cells.Item(30, 7).Value = "'21"
$ws.Cells.Item(31, 7).Value = "'21"
$ws.Cells.Item(32, 7).Value = "'21"
$ws.Cells.Item(33, 7).Value = "'21"
$ws.Cells.Item(34, 7).Value = "'21"
$ws.Cells.Item(35, 7).Value = "'21"
$ws.Cells.Item(36, 7).Value = "'21"
$ws.Cells.Item(37, 7).Value = "'21"
$ws.Cells.Item(38, 7).Value = "'21"
$ws.Cells.Item(39, 7).Value = "'21"
$ws.Cells.Item(40, 4).Value = "'0.04665"
$ws.Cells.Item(40, 7).Value = "'21"
$ws.Cells.Item(41, 4).Value = "'0.006883"
$ws.Cells.Item(41, 7).Value = "'21"
$ws.Cells.Item(42, 4).Value = "'0.1155"
$ws.Cells.Item(42, 7).Value = "'21"
$ws.Cells.Item(43, 4).Value = "'0.003354"
$ws.Cells.Item(43, 7).Value = "'21"
$ws.Cells.Item(44, 4).Value = "'0.01226"
$ws.Cells.Item(44, 7).Value = "'21"
$ws.Cells.Item(45, 4).Value = "'0.00006261"
$ws.Cells.Item(45, 7).Value = "'21"
$ws.Cells.Item(46, 4).Value = "'0.00000000751"
$ws.Cells.Item(46, 7).Value = "'21"
$ws.Cells.Item(47, 4).Value = "'0.7902"
$ws.Cells.Item(47, 7).Value = "'21"
$ws.Cells.Item(48, 4).Value = "'0.1602"
$ws.Cells.Item(48, 7).Value = "'21"
$ws.Cells.Item(49, 4).Value = "'0.00002104"
$ws.Cells.Item(49, 7).Value = "'21"
$ws.Cells.Item(50, 4).Value = "'0.01242"
$ws.Cells.Item(50, 7).Value = "'21"
$ws.Cells.Item(51, 7).Value = "'21"
